# Update loading_percent values for case with 380 kV.
# Columns B, D, E, F, G, H, I, N are updated for rows 2-25 (24 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B", "D", "E", "F", "G", "H", "I", "N")

$data = @(
    @(7.973936579530329, 4.4530821209041, 16.80072893514019, 31.66688530865378, 47.97825748599407, 14.99599398427061, 17.51723746316919, 17.13588176064243),
    @(7.726105456002573, 4.256110385767948, 15.81414347791596, 30.1364997654836, 45.07585997036803, 14.56739984798534, 17.43545559361282, 16.94880368210766),
    @(7.569968192699863, 4.129921971931488, 15.18345797957187, 29.16934848132706, 43.21230408299177, 14.30402308097859, 17.38881070702378, 16.8348163882461),
    @(7.505429865113701, 4.077210539987413, 14.92044224378858, 28.76891537651177, 42.43304870019238, 14.19681494693921, 17.37072111759998, 16.78863424108457),
    @(7.494660968687235, 4.06838101841469, 14.87641457860602, 28.70206186205763, 42.30247733761876, 14.17902561583299, 17.36777342702464, 16.78098333545054),
    @(7.569101375573601, 4.129216260298202, 15.17993479371136, 29.16397279584906, 43.20187414878671, 14.30257651222966, 17.38856299824779, 16.83419240957982),
    @(7.889357710724763, 4.386272154184288, 16.46586561567118, 31.14524790463039, 46.99476461840648, 14.84836263646827, 17.48831115075836, 17.07122325684947),
    @(8.48245444634245, 4.847626781533176, 18.96963095844109, 34.79019700568814, 53.76441982954865, 15.90964749827232, 17.71126666152004, 17.54092831277462),
    @(8.892805469389883, 5.159375015272185, 20.7005264712871, 37.29712272872882, 58.30908890538318, 16.67484774318615, 17.89037526883937, 17.88621396499388),
    @(9.073239364152251, 5.295119806458914, 21.44630217105968, 38.39657939506925, 60.28053908437305, 17.01811053947613, 17.97487022382017, 18.04275315944065),
    @(9.140615808877813, 5.345639535607845, 21.72277565408117, 38.8067745188044, 61.01311911472015, 17.14727530363974, 18.00727319226982, 18.10190618682361),
    @(9.126148118074395, 5.334798663661699, 21.66349534754383, 38.71870855538847, 60.85596817346871, 17.11949578527255, 18.00027699228534, 18.08917283935257),
    @(9.078801786399699, 5.299293861980987, 21.46916648747204, 38.43045085574376, 60.34108956117144, 17.02875414856346, 17.97752804593549, 18.04762249193302),
    @(9.049675604494446, 5.277430761849357, 21.34936291773622, 38.25307759903486, 60.02388868579952, 16.97306174473784, 17.96364578513084, 18.02215401479311),
    @(8.880883512103347, 5.150380618130106, 20.65095500525945, 37.2244235607972, 58.17830377378787, 16.65230639989389, 17.88491171931956, 17.87596887941638),
    @(8.77569792055426, 5.070875926456215, 20.21187991088275, 36.58269053228319, 57.02138715638986, 16.45420217812308, 17.83736449730068, 17.78611779882747),
    @(8.71461368730956, 5.024575959885126, 19.95541536592146, 36.20973971098365, 56.34693520985559, 16.33981001259289, 17.8103022940806, 17.73439049486328),
    @(8.693832984641853, 5.008801897977272, 19.86790610021658, 36.08281351125076, 56.11703346922739, 16.30100592805423, 17.80118939117584, 17.71686987832945),
    @(8.786955979029544, 5.079398517153265, 20.25902557245064, 36.65140384155789, 57.14547825128386, 16.47533801822657, 17.842396607743, 17.79568785271756),
    @(9.092734726993198, 5.309746544556469, 21.52640628333927, 38.51528771957199, 60.49270200548502, 17.0554304110521, 17.9841991494933, 18.05983060728496),
    @(9.287021218569393, 5.455137147564349, 22.32013858859541, 39.69753703054367, 62.59886043891643, 17.42971607441627, 18.07923203987227, 18.23171317604657),
    @(9.183851410766975, 5.378014062561288, 21.89965602861509, 39.06990661277872, 61.4822596397713, 17.23043429362993, 18.02830472145925, 18.14006050298211),
    @(8.781868111808391, 5.075547295935355, 20.23772358888641, 36.62035102466096, 57.08940569978319, 16.46578405071612, 17.84012073718834, 17.79136144795664),
    @(8.326196084739871, 4.727527600169029, 18.29498466018679, 33.83266512629792, 52.00717024423898, 15.62447975441807, 17.64816122966614, 17.41362653806358)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $columns.Length; $j++) {
        $col = $columns[$j]
        $ws.Range("$col$row").Value = $rowValues[$j]
    }
}
